$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table (B4:G8) gains a new "Deallocation" row, following the same
# layout as the existing rows (e.g. row 8 "Dereplikation"): a label in
# column B and "a" markers across C:G, using the same formatting.
$ws.Range("B8:G8").Copy()
$ws.Range("B9:G9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("B9").Value = "Deallocation"
$ws.Range("C9:G9").Value = "a"
$ws.Rows(9).RowHeight = 16.5

$ws.Range("J10").Select()
